$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 12120533
$ws.Range("I19").Value = 10435505
$ws.Range("K19").Value = 10435505
$ws.Range("M19").Value = -10435330
$ws.Range("H40").Value = 1755.2354
$ws.Range("I40").Value = 1533.2858
$ws.Range("J40").Value = 2113.7693
$ws.Range("K40").Value = 1533.2858
$ws.Range("L40").Value = 2113.7693
$ws.Range("M40").Value = -1358.2858
$ws.Range("N40").Value = -2463.7693
$ws.Range("H92").Value = 1182
$ws.Range("I92").Value = 385
$ws.Range("J92").Value = 4370
$ws.Range("K92").Value = 385
$ws.Range("L92").Value = 4370
$ws.Range("M92").Value = 863
$ws.Range("N92").Value = -6866
$ws.Range("H99").Value = 992
$ws.Range("I99").Value = 476
$ws.Range("J99").Value = 1443.5
$ws.Range("K99").Value = 1428
$ws.Range("L99").Value = 4330.5
$ws.Range("M99").Value = 70
$ws.Range("N99").Value = -7326.5
$ws.Range("H112").Value = 1931
$ws.Range("I112").Value = 199.5
$ws.Range("J112").Value = 2197.3845
$ws.Range("K112").Value = 598.5
$ws.Range("L112").Value = 6592.1535
$ws.Range("M112").Value = 509.5
$ws.Range("N112").Value = -8808.1535
$ws.Range("H113").Value = 1972.1111
$ws.Range("I113").Value = 1750
$ws.Range("J113").Value = 1999.875
$ws.Range("K113").Value = 1750
$ws.Range("L113").Value = 1999.875
$ws.Range("M113").Value = 1504
$ws.Range("N113").Value = -8507.875
$ws.Range("H116").Value = 1751038.1
$ws.Range("J116").Value = 3383.5417
$ws.Range("L116").Value = 3383.5417
$ws.Range("N116").Value = -10267.5417
$ws.Range("H135").Value = 274.63635
$ws.Range("I135").Value = 274.63635
$ws.Range("K135").Value = 2471.72715
$ws.Range("M135").Value = 63.27285000000029

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19235810
$ws.Range("I32").Value = 23811794
$ws.Range("K32").Value = 23811794
$ws.Range("M32").Value = -23811507
$ws.Range("H97").Value = 406.82858
$ws.Range("I97").Value = 375.2647
$ws.Range("J97").Value = 1480
$ws.Range("K97").Value = 375.2647
$ws.Range("L97").Value = 1480
$ws.Range("M97").Value = 120.7353
$ws.Range("N97").Value = -2472
$ws.Range("H122").Value = 1482943.1
$ws.Range("I122").Value = 1544607.4
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4633822.199999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4631372.199999999
$ws.Range("N122").Value = -13900

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 510.64706
$ws.Range("I64").Value = 648.75
$ws.Range("K64").Value = 648.75
$ws.Range("M64").Value = -423.75
$ws.Range("H67").Value = 510.64706
$ws.Range("I67").Value = 648.75
$ws.Range("K67").Value = 648.75
$ws.Range("M67").Value = 131.25
$ws.Range("H80").Value = 130.875
$ws.Range("I80").Value = 161.5
$ws.Range("J80").Value = 100.25
$ws.Range("K80").Value = 161.5
$ws.Range("L80").Value = 100.25
$ws.Range("M80").Value = 836.5
$ws.Range("N80").Value = -2096.25
$ws.Range("H83").Value = 130.875
$ws.Range("I83").Value = 161.5
$ws.Range("J83").Value = 100.25
$ws.Range("K83").Value = 807.5
$ws.Range("L83").Value = 501.25
$ws.Range("M83").Value = 4184.5
$ws.Range("N83").Value = -10485.25
$ws.Range("H86").Value = 3110.75
$ws.Range("I86").Value = 2944
$ws.Range("K86").Value = 2944
$ws.Range("M86").Value = -1821
$ws.Range("H89").Value = 3110.75
$ws.Range("I89").Value = 2944
$ws.Range("K89").Value = 14720
$ws.Range("M89").Value = -9104

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2523.5344
$ws.Range("I31").Value = 1723.591
$ws.Range("K31").Value = 1723.591
$ws.Range("M31").Value = -1428.591
$ws.Range("H34").Value = 2523.5344
$ws.Range("I34").Value = 1723.591
$ws.Range("K34").Value = 1723.591
$ws.Range("M34").Value = -1521.591
$ws.Range("H107").Value = 514.1739
$ws.Range("I107").Value = 293.33334
$ws.Range("J107").Value = 928.25
$ws.Range("K107").Value = 293.33334
$ws.Range("L107").Value = 928.25
$ws.Range("M107").Value = 1626.66666
$ws.Range("N107").Value = -4768.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 833623.0600000001
$ws.Range("J92").Value = 338.625
$ws.Range("L92").Value = 1015.875
$ws.Range("N92").Value = -3511.875
$ws.Range("H97").Value = 949.3570999999999
$ws.Range("J97").Value = 1400
$ws.Range("L97").Value = 4200
$ws.Range("N97").Value = -5192
$ws.Range("H131").Value = 811.2969000000001
$ws.Range("I131").Value = 398.0625
$ws.Range("J131").Value = 949.0417
$ws.Range("K131").Value = 1194.1875
$ws.Range("L131").Value = 2847.1251
$ws.Range("M131").Value = 3845.8125
$ws.Range("N131").Value = -12927.1251

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 301
$ws.Range("I107").Value = 244.92308
$ws.Range("K107").Value = 244.92308
$ws.Range("M107").Value = 1675.07692
$ws.Range("H113").Value = 13605.556
$ws.Range("I113").Value = 2064.2856
$ws.Range("J113").Value = 54000
$ws.Range("K113").Value = 2064.2856
$ws.Range("L113").Value = 54000
$ws.Range("M113").Value = 105.7143999999998
$ws.Range("N113").Value = -58340
$ws.Range("H126").Value = 2502.4
$ws.Range("I126").Value = 1442.8
$ws.Range("K126").Value = 4328.4
$ws.Range("M126").Value = -1858.4
$ws.Range("H132").Value = 3689.5881
$ws.Range("I132").Value = 3152.375
$ws.Range("K132").Value = 9457.125
$ws.Range("M132").Value = -6927.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 112222920
$ws.Range("I22").Value = 168333740
$ws.Range("J22").Value = 1265
$ws.Range("K22").Value = 168333740
$ws.Range("L22").Value = 1265
$ws.Range("M22").Value = -168333445
$ws.Range("N22").Value = -1855
$ws.Range("H27").Value = 112222920
$ws.Range("I27").Value = 168333740
$ws.Range("J27").Value = 1265
$ws.Range("K27").Value = 168333740
$ws.Range("L27").Value = 1265
$ws.Range("M27").Value = -168333633
$ws.Range("N27").Value = -1479
$ws.Range("H61").Value = 3894.3428
$ws.Range("I61").Value = 4768.2856
$ws.Range("J61").Value = 398.57144
$ws.Range("K61").Value = 4768.2856
$ws.Range("L61").Value = 398.57144
$ws.Range("M61").Value = -4566.2856
$ws.Range("N61").Value = -802.5714399999999
$ws.Range("H113").Value = 3894.3428
$ws.Range("I113").Value = 4768.2856
$ws.Range("J113").Value = 398.57144
$ws.Range("K113").Value = 4768.2856
$ws.Range("L113").Value = 398.57144
$ws.Range("M113").Value = -2598.2856
$ws.Range("N113").Value = -4738.57144
$ws.Range("H132").Value = 4166.75
$ws.Range("I132").Value = 3687.875
$ws.Range("J132").Value = 5124.5
$ws.Range("K132").Value = 11063.625
$ws.Range("L132").Value = 15373.5
$ws.Range("M132").Value = -8533.625
$ws.Range("N132").Value = -20433.5
$ws.Range("H133").Value = 49066.3
$ws.Range("J133").Value = 49066.3
$ws.Range("L133").Value = 49066.3
$ws.Range("N133").Value = -54126.3

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1019.46155
$ws.Range("J100").Value = 1113.909
$ws.Range("L100").Value = 2227.818
$ws.Range("N100").Value = -3309.818
$ws.Range("H113").Value = 440.8125
$ws.Range("I113").Value = 311.76923
$ws.Range("K113").Value = 935.30769
$ws.Range("M113").Value = 1234.69231
